$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(4, 7, 14, 195, 179, 3, 1, 0, 91),
    @(10, 8, 3, 234, 230, 3, 1, 1, 8),
    @(12, 100, 14, 198, 184, 3, 1, 1, 1100)
)

$startRow = 21
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
